$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '37.729.79'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '2.078.25'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''58.07'
$ws.Range("E8").Value = '  -1.48%  '
$ws.Range("D10").Value = '''0.0781'
$ws.Range("E10").Value = '  -1.05%  '
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").Value = '2.385.31'
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("D17").Value = '2.081.25'
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("D18").Value = '37.662.15'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").Value = '''70.26'
$ws.Range("E20").Value = '  -2.20%  '
$ws.Range("D21").Value = '0.0₃0831'
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("D22").Value = '''227.67'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '''2.39'
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D26").Value = '''9.94'
$ws.Range("E26").Value = '  +3.33%  '
$ws.Range("D27").Value = '''169.56'
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("E28").Value = '  -3.90%  '
$ws.Range("D29").Value = '''19.37'
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("E30").Value = '  -3.16%  '
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").Value = '''4.62'
$ws.Range("E32").Value = '  -2.70%  '
$ws.Range("D33").Value = '''0.0631'
$ws.Range("E33").Value = '  -0.63%  '
$ws.Range("D34").Value = '''4.66'
$ws.Range("E34").Value = '  -0.44%  '
$ws.Range("D35").Value = '''2.53'
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("E36").Value = '  -0.42%  '
$ws.Range("D37").Value = '''3.33'
$ws.Range("E37").Value = '  -3.28%  '
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("E39").Value = '  -1.64%  '
$ws.Range("D40").Value = '''0.0228'
$ws.Range("E40").Value = '  +3.83%  '
$ws.Range("D41").Value = '''98.46'
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").Value = '''2.92'
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").Value = '''0.0958'
$ws.Range("E43").Value = '  -2.36%  '
$ws.Range("D44").Value = '1.489.58'
$ws.Range("E44").Value = '  +2.55%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '''17.03'
$ws.Range("E45").Value = '  -2.34%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = '''1.19'
$ws.Range("E46").Value = '  +3.11%  '
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("E48").Value = '  -1.70%  '
$ws.Range("D50").Value = '''2.96'
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("D51").Value = '2.269.55'
$ws.Range("E51").Value = '  -0.31%  '
